# soberano responsibilities test-cases matrix:
# add a new test-case row for "user21" (product category check recording
# test) right above the existing workshop1Worker/workshop2Worker rows,
# pushing everything from row 34 down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 34 - this shifts the old rows 34-38 down to 35-39
# and (matching Excel's native behaviour) copies the formatting of the row
# above (33, the yellow "new test case" highlight) into the new row.
$ws.Rows("34:34").Insert()

# New test case: user21 is only assigned the SystemAdmin (column N) role.
$ws.Range("A34").Value = "user21"
$ws.Range("B34").Value = "not assigned"
$ws.Range("C34").Value = "not assigned"
$ws.Range("D34").Value = "not assigned"
$ws.Range("E34").Value = "not assigned"
$ws.Range("F34").Value = "not assigned"
$ws.Range("G34").Value = "not assigned"
$ws.Range("H34").Value = "not assigned"
$ws.Range("I34").Value = "not assigned"
$ws.Range("J34").Value = "not assigned"
$ws.Range("K34").Value = "not assigned"
$ws.Range("L34").Value = "not assigned"
$ws.Range("M34").Value = "not assigned"
$ws.Range("N34").Value = "assigned"

# Match the row heights used throughout the sheet: the new highlighted row
# uses the same compact height as the other highlighted test-case rows,
# while the rows that got pushed down keep their original heights (the
# Insert() above reset them).
$ws.Rows("34:34").RowHeight = 13.8
$ws.Rows("35:35").RowHeight = 15
$ws.Rows("36:36").RowHeight = 15
$ws.Rows("37:37").RowHeight = 15
$ws.Rows("39:39").RowHeight = 14.9

# Leave the final selection on the newly recorded cell.
$ws.Range("N34").Select() | Out-Null
